# Update DummyData: rotate the dummy "hashed_password" value used for every
# user row, and update the active sheet's selection to reflect the password
# column (D2:D12) now being highlighted instead of just the last cell (D12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")
$ws.Activate()

# All rows share the same dummy bcrypt hash in column D (rows 2-12).
# Re-write them all to the new hash value.
$newHash = '$2a$12$wAjivtgOA8ghj4wU6.P1weHsvGfLWfWaeaW3RKrST5gt72smmd5hW'
$ws.Range("D2:D12").Value = $newHash

# Reflect the updated selection: whole password column body selected,
# active cell on the first data row.
$ws.Range("D2:D12").Select()
